# Append a new time-logging entry (row 10) to the sheet, mirroring the
# structure of the existing rows (3-9): a date, a start time, an end
# time, a duration formula (End - Start) and a category label.
#
# The source rows use custom number formats (date / time-of-day) that
# live in styles.xml, so we copy the formatting from an existing sibling
# cell (rather than writing NumberFormat strings, which would mint new
# duplicate style entries) and only then overwrite the copied value with
# the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 10
$prevRow = 9

# Carry over formatting (date style, time style, text style) from the
# previous row's cells.
$ws.Range("B$prevRow").Copy($ws.Range("B$row"))
$ws.Range("C$prevRow").Copy($ws.Range("C$row"))
$ws.Range("D$prevRow").Copy($ws.Range("D$row"))
$ws.Range("F$prevRow").Copy($ws.Range("F$row"))
# E3 is the head of the shared "D-C" duration formula; copy it so the new
# cell inherits the same (time-of-day) number format.
$ws.Range("E3").Copy($ws.Range("E$row"))

# New entry: 2025-12-02, 12:15 - 14:30, NanoGPT.
$ws.Cells.Item($row, 2).Value = 45993
$ws.Cells.Item($row, 3).Value = 0.51041666666666663
$ws.Cells.Item($row, 4).Value = 0.60416666666666663
$ws.Cells.Item($row, 5).Formula = "=D$row-C$row"
$ws.Cells.Item($row, 6).Value = "NanoGPT"
